$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row content swaps (columns B:AC) for paired/grouped fixtures ----
# Each cycle below rotates the B:AC row content among the listed rows,
# i.e. row r takes on the former B:AC values of row mapping[r], while
# column A (the running match index) stays untouched on every row.

# Cycle: [8, 9]
$v8 = $ws.Range("B8:AC8").Value2
$v9 = $ws.Range("B9:AC9").Value2
$ws.Range("B8:AC8").Value2 = $v9
$ws.Range("B9:AC9").Value2 = $v8

# Cycle: [23, 24]
$v23 = $ws.Range("B23:AC23").Value2
$v24 = $ws.Range("B24:AC24").Value2
$ws.Range("B23:AC23").Value2 = $v24
$ws.Range("B24:AC24").Value2 = $v23

# Cycle: [29, 31, 30]
$v29 = $ws.Range("B29:AC29").Value2
$v31 = $ws.Range("B31:AC31").Value2
$v30 = $ws.Range("B30:AC30").Value2
$ws.Range("B29:AC29").Value2 = $v31
$ws.Range("B31:AC31").Value2 = $v30
$ws.Range("B30:AC30").Value2 = $v29

# Cycle: [32, 33]
$v32 = $ws.Range("B32:AC32").Value2
$v33 = $ws.Range("B33:AC33").Value2
$ws.Range("B32:AC32").Value2 = $v33
$ws.Range("B33:AC33").Value2 = $v32

# Cycle: [39, 41]
$v39 = $ws.Range("B39:AC39").Value2
$v41 = $ws.Range("B41:AC41").Value2
$ws.Range("B39:AC39").Value2 = $v41
$ws.Range("B41:AC41").Value2 = $v39

# Cycle: [63, 64]
$v63 = $ws.Range("B63:AC63").Value2
$v64 = $ws.Range("B64:AC64").Value2
$ws.Range("B63:AC63").Value2 = $v64
$ws.Range("B64:AC64").Value2 = $v63

# Cycle: [89, 90]
$v89 = $ws.Range("B89:AC89").Value2
$v90 = $ws.Range("B90:AC90").Value2
$ws.Range("B89:AC89").Value2 = $v90
$ws.Range("B90:AC90").Value2 = $v89

# Cycle: [105, 106]
$v105 = $ws.Range("B105:AC105").Value2
$v106 = $ws.Range("B106:AC106").Value2
$ws.Range("B105:AC105").Value2 = $v106
$ws.Range("B106:AC106").Value2 = $v105

# Cycle: [129, 130]
$v129 = $ws.Range("B129:AC129").Value2
$v130 = $ws.Range("B130:AC130").Value2
$ws.Range("B129:AC129").Value2 = $v130
$ws.Range("B130:AC130").Value2 = $v129

# Cycle: [135, 141, 140, 138, 139, 137]
$v135 = $ws.Range("B135:AC135").Value2
$v141 = $ws.Range("B141:AC141").Value2
$v140 = $ws.Range("B140:AC140").Value2
$v138 = $ws.Range("B138:AC138").Value2
$v139 = $ws.Range("B139:AC139").Value2
$v137 = $ws.Range("B137:AC137").Value2
$ws.Range("B135:AC135").Value2 = $v141
$ws.Range("B141:AC141").Value2 = $v140
$ws.Range("B140:AC140").Value2 = $v138
$ws.Range("B138:AC138").Value2 = $v139
$ws.Range("B139:AC139").Value2 = $v137
$ws.Range("B137:AC137").Value2 = $v135

# Cycle: [172, 174]
$v172 = $ws.Range("B172:AC172").Value2
$v174 = $ws.Range("B174:AC174").Value2
$ws.Range("B172:AC172").Value2 = $v174
$ws.Range("B174:AC174").Value2 = $v172

# Cycle: [191, 192]
$v191 = $ws.Range("B191:AC191").Value2
$v192 = $ws.Range("B192:AC192").Value2
$ws.Range("B191:AC191").Value2 = $v192
$ws.Range("B192:AC192").Value2 = $v191

# Cycle: [213, 214]
$v213 = $ws.Range("B213:AC213").Value2
$v214 = $ws.Range("B214:AC214").Value2
$ws.Range("B213:AC213").Value2 = $v214
$ws.Range("B214:AC214").Value2 = $v213

# ---- Direct odds updates for upcoming fixtures (rows 218-221) ----
# Row 218
$ws.Range("O218").Value2 = 3.2
$ws.Range("P218").Value2 = 1.909
$ws.Range("R218").Value2 = 1.875
$ws.Range("S218").Value2 = 1.925
$ws.Range("U218").Value2 = 2
$ws.Range("V218").Value2 = 1.8

# Row 219
$ws.Range("N219").Value2 = 2.05
$ws.Range("O219").Value2 = 3.25
$ws.Range("P219").Value2 = 3.1
$ws.Range("R219").Value2 = 1.875
$ws.Range("S219").Value2 = 1.925
$ws.Range("U219").Value2 = 1.9
$ws.Range("V219").Value2 = 1.9

# Row 220
$ws.Range("N220").Value2 = 1.615
$ws.Range("O220").Value2 = 3.5
$ws.Range("P220").Value2 = 4.75
$ws.Range("Q220").Value2 = -0.75
$ws.Range("R220").Value2 = 1.8
$ws.Range("S220").Value2 = 2
$ws.Range("U220").Value2 = 2.025
$ws.Range("V220").Value2 = 1.775

# Row 221
$ws.Range("N221").Value2 = 1.85
$ws.Range("P221").Value2 = 4.2
$ws.Range("R221").Value2 = 1.95
$ws.Range("S221").Value2 = 1.85
$ws.Range("U221").Value2 = 1.825
$ws.Range("V221").Value2 = 1.975

